# Recon_PaidInFull_Jan2026.xlsx update: refresh Paid-in-Full reconciliation tape rows 4-15
# with new loan IDs/investors/loan types/payoff data for the regenerated January 2026 tape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "MSR100034"
$ws.Range("B4").Value = "FHLMC"
$ws.Range("C4").Value = "FHA"
$ws.Range("D4").Value = 46031
$ws.Range("E4").Value = 485000
$ws.Range("F4").Value = 424666.91
$ws.Range("G4").Value = 426984.98
$ws.Range("H4").Value = 2229.5
$ws.Range("I4").Value = 88.56999999999999
$ws.Range("J4").Value = 426984.98
$ws.Range("K4").Value = 0.063
$ws.Range("L4").Value = "Full Payoff"
$ws.Range("M4").Value = "Removed from portfolio 01/09/2026"
$ws.Range("A5").Value = "MSR100100"
$ws.Range("B5").Value = "FNMA"
$ws.Range("C5").Value = "FHA"
$ws.Range("D5").Value = 46049
$ws.Range("E5").Value = 390000
$ws.Range("F5").Value = 291530.86
$ws.Range("G5").Value = 292526.74
$ws.Range("H5").Value = 760.41
$ws.Range("I5").Value = 235.47
$ws.Range("J5").Value = 292526.74
$ws.Range("K5").Value = 0.0313
$ws.Range("L5").Value = "Refinance Payoff"
$ws.Range("M5").Value = "Removed from portfolio 01/27/2026"
$ws.Range("A6").Value = "MSR100241"
$ws.Range("B6").Value = "GNMA"
$ws.Range("C6").Value = "FHA"
$ws.Range("D6").Value = 46038
$ws.Range("E6").Value = 300000
$ws.Range("F6").Value = 265375.94
$ws.Range("G6").Value = 266391.6
$ws.Range("H6").Value = 683.34
$ws.Range("I6").Value = 332.32
$ws.Range("J6").Value = 266391.6
$ws.Range("K6").Value = 0.0309
$ws.Range("L6").Value = "Refinance Payoff"
$ws.Range("M6").Value = "Removed from portfolio 01/16/2026"
$ws.Range("A7").Value = "MSR100250"
$ws.Range("B7").Value = "FNMA"
$ws.Range("C7").Value = "FHA"
$ws.Range("D7").Value = 46032
$ws.Range("E7").Value = 305000
$ws.Range("F7").Value = 106147.08
$ws.Range("G7").Value = 106570.37
$ws.Range("H7").Value = 397.17
$ws.Range("I7").Value = 26.12
$ws.Range("J7").Value = 106570.37
$ws.Range("K7").Value = 0.0449
$ws.Range("L7").Value = "Refinance Payoff"
$ws.Range("M7").Value = "Removed from portfolio 01/10/2026"
$ws.Range("A8").Value = "MSR100252"
$ws.Range("B8").Value = "FNMA"
$ws.Range("C8").Value = "Conventional"
$ws.Range("D8").Value = 46050
$ws.Range("E8").Value = 470000
$ws.Range("F8").Value = 365421
$ws.Range("G8").Value = 366849.04
$ws.Range("H8").Value = 1367.28
$ws.Range("I8").Value = 60.76
$ws.Range("J8").Value = 366849.04
$ws.Range("K8").Value = 0.0449
$ws.Range("L8").Value = "Refinance Payoff"
$ws.Range("M8").Value = "Removed from portfolio 01/28/2026"
$ws.Range("A9").Value = "MSR100253"
$ws.Range("B9").Value = "FHLMC"
$ws.Range("C9").Value = "VA"
$ws.Range("D9").Value = 46036
$ws.Range("E9").Value = 355000
$ws.Range("F9").Value = 304042.36
$ws.Range("G9").Value = 305430.66
$ws.Range("H9").Value = 1195.9
$ws.Range("I9").Value = 192.4
$ws.Range("J9").Value = 305430.66
$ws.Range("K9").Value = 0.0472
$ws.Range("L9").Value = "Refinance Payoff"
$ws.Range("M9").Value = "Removed from portfolio 01/14/2026"
$ws.Range("A10").Value = "MSR100264"
$ws.Range("B10").Value = "GNMA"
$ws.Range("C10").Value = "USDA"
$ws.Range("D10").Value = 46043
$ws.Range("E10").Value = 545000
$ws.Range("F10").Value = 502226.21
$ws.Range("G10").Value = 505073.62
$ws.Range("H10").Value = 2599.02
$ws.Range("I10").Value = 248.39
$ws.Range("J10").Value = 505073.62
$ws.Range("K10").Value = 0.0621
$ws.Range("L10").Value = "Refinance Payoff"
$ws.Range("M10").Value = "Removed from portfolio 01/21/2026"
$ws.Range("A11").Value = "MSR100625"
$ws.Range("B11").Value = "FNMA"
$ws.Range("C11").Value = "Conventional"
$ws.Range("D11").Value = 46040
$ws.Range("E11").Value = 515000
$ws.Range("F11").Value = 450060.9
$ws.Range("G11").Value = 451942.7
$ws.Range("H11").Value = 1796.49
$ws.Range("I11").Value = 85.31
$ws.Range("J11").Value = 451942.7
$ws.Range("K11").Value = 0.0479
$ws.Range("L11").Value = "Refinance Payoff"
$ws.Range("M11").Value = "Removed from portfolio 01/18/2026"
$ws.Range("A12").Value = "MSR100709"
$ws.Range("B12").Value = "FNMA"
$ws.Range("C12").Value = "FHA"
$ws.Range("D12").Value = 46031
$ws.Range("E12").Value = 390000
$ws.Range("F12").Value = 311889.86
$ws.Range("G12").Value = 313016.71
$ws.Range("H12").Value = 836.9
$ws.Range("I12").Value = 289.95
$ws.Range("J12").Value = 313016.71
$ws.Range("K12").Value = 0.0322
$ws.Range("L12").Value = "Full Payoff"
$ws.Range("M12").Value = "Removed from portfolio 01/09/2026"
$ws.Range("A13").Value = "MSR100911"
$ws.Range("B13").Value = "FNMA"
$ws.Range("C13").Value = "FHA"
$ws.Range("D13").Value = 46037
$ws.Range("E13").Value = 250000
$ws.Range("F13").Value = 212771.66
$ws.Range("G13").Value = 213785.31
$ws.Range("H13").Value = 790.8
$ws.Range("I13").Value = 222.85
$ws.Range("J13").Value = 213785.31
$ws.Range("K13").Value = 0.0446
$ws.Range("L13").Value = "Refinance Payoff"
$ws.Range("M13").Value = "Removed from portfolio 01/15/2026"
$ws.Range("A14").Value = "MSR100913"
$ws.Range("B14").Value = "FHLMC"
$ws.Range("C14").Value = "Conventional"
$ws.Range("D14").Value = 46033
$ws.Range("E14").Value = 425000
$ws.Range("F14").Value = 374742.05
$ws.Range("G14").Value = 375838.44
$ws.Range("H14").Value = 949.35
$ws.Range("I14").Value = 147.04
$ws.Range("J14").Value = 375838.44
$ws.Range("K14").Value = 0.0304
$ws.Range("L14").Value = "Sale Payoff"
$ws.Range("M14").Value = "Removed from portfolio 01/11/2026"
$ws.Range("A15").Value = "MSR100952"
$ws.Range("B15").Value = "FHLMC"
$ws.Range("C15").Value = "FHA"
$ws.Range("D15").Value = 46052
$ws.Range("E15").Value = 335000
$ws.Range("F15").Value = 233578.42
$ws.Range("G15").Value = 234850.72
$ws.Range("H15").Value = 924.58
$ws.Range("I15").Value = 347.72
$ws.Range("J15").Value = 234850.72
$ws.Range("K15").Value = 0.0475
$ws.Range("L15").Value = "Full Payoff"
$ws.Range("M15").Value = "Removed from portfolio 01/30/2026"
